# Refresh cryptocurrency price / 1h-volume-change data, and fix the
# position of two row-pairs (Cosmos/Toncoin, FirstDigitalUSD/Filecoin)
# whose rank order changed in the upstream feed.
#
# Price ("D") cells are plain numeric-looking strings in the source feed
# (e.g. "0.999", "307.61") but are stored as TEXT in the workbook, not
# numbers (the existing sheet uses inline strings for every data cell).
# Assigning such a string straight to .Value would make Excel
# auto-convert it to a real number, so for any new price value that is
# unambiguously numeric we first flip the cell to Text format ("@") -
# this mirrors what happens when a user / feed writes a text price into
# a plain cell and matches the original columns text typing.
# Values that already contain a second "." (thousand separator, e.g.
# "42.086.02") or other non-numeric characters stay text on their own,
# so no extra formatting call is needed for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.086.02"
$ws.Range("E2").Value = "  -0.90%  "

# Row 3
$ws.Range("D3").Value = "2.260.31"
$ws.Range("E3").Value = "  -1.18%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.61"
$ws.Range("E5").Value = "  +0.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.05"
$ws.Range("E6").Value = "  +0.23%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.525"
$ws.Range("E7").Value = "  -1.34%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -1.49%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.15"
$ws.Range("E10").Value = "  -3.37%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").Value = "  -2.37%  "

# Row 12
$ws.Range("E12").Value = "  +0.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.79"
$ws.Range("E13").Value = "  +1.09%  "

# Row 14
$ws.Range("D14").Value = "2.608.26"
$ws.Range("E14").Value = "  -1.20%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.55"
$ws.Range("E15").Value = "  -0.73%  "

# Row 16
$ws.Range("D16").Value = "2.260.86"
$ws.Range("E16").Value = "  -1.51%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.787"
$ws.Range("E17").Value = "  -2.26%  "

# Row 18
$ws.Range("D18").Value = "41.895.31"
$ws.Range("E18").Value = "  -1.10%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.19"
$ws.Range("E19").Value = "  -5.47%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0898"
$ws.Range("E20").Value = "  -2.36%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.94"
$ws.Range("E21").Value = "  -1.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.43"
$ws.Range("E22").Value = "  -0.59%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.84"
$ws.Range("E23").Value = "  -2.87%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.58"
$ws.Range("E24").Value = "  -0.97%  "

# Row 25
$ws.Range("E25").Value = "  +0.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.48"
$ws.Range("E27").Value = "  -2.18%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.36"
$ws.Range("E28").Value = "  -0.76%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("E29").Value = "  +0.54%  "

# Row 30
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.49"
$ws.Range("E30").Value = "  -1.19%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "163.83"
$ws.Range("E31").Value = "  +1.55%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.21"
$ws.Range("E32").Value = "  -2.29%  "

# Row 33
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.03%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("E34").Value = "  +0.27%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0733"
$ws.Range("E35").Value = "  -2.56%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.28"
$ws.Range("E36").Value = "  -0.78%  "

# Row 37
$ws.Range("E37").Value = "  +0.05%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.104"
$ws.Range("E38").Value = "  -4.16%  "

# Row 39
$ws.Range("E39").Value = "  -1.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.81"
$ws.Range("E40").Value = "  -3.83%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.13"
$ws.Range("E41").Value = "  -1.32%  "

# Row 42
$ws.Range("E42").Value = "  -5.65%  "

# Row 43
$ws.Range("D43").Value = "1.950.56"
$ws.Range("E43").Value = "  -2.59%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.78"
$ws.Range("E44").Value = "  -3.00%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0280"
$ws.Range("E45").Value = "  -2.09%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.91"
$ws.Range("E46").Value = "  -3.61%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.77"
$ws.Range("E47").Value = "  -4.97%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.24"
$ws.Range("E48").Value = "  -1.51%  "

# Row 49
$ws.Range("D49").Value = "2.480.76"
$ws.Range("E49").Value = "  -1.15%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.02"
$ws.Range("E50").Value = "  +0.02%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.28"
$ws.Range("E51").Value = "  -1.86%  "
